# Update daily power records
# - complete the in-progress record in row 25 (End Time)
# - add a new complete record in row 26
# - add a new (currently blank) record placeholder in row 27
# - extend the Excel table / autofilter / dimension to cover the new rows
# - update the active selection to reflect where the user ended up working

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 25: fill in the missing "End Time" (C25) for the existing record ---
$ws.Range("C25").Value = 0.56805555555555554

# --- Duplicate row 25 (now complete, with correct formatting/formulas) down
#     into rows 26 and 27 so the new rows inherit the same number formats,
#     styles and calculated-column formulas used throughout the table. ---
$ws.Rows("25:25").Copy()
$ws.Rows("26:26").Insert(-4121)   # -4121 = xlShiftDown

$ws.Rows("25:25").Copy()
$ws.Rows("27:27").Insert(-4121)   # -4121 = xlShiftDown

# --- Row 26: a brand new, complete record ---
$ws.Range("A26").Value = 43352
$ws.Range("B26").Value = 0.57500000000000007
$ws.Range("C26").Value = 0.72361111111111109

# --- Row 27: new record with Date/Start/End not entered yet, so clear those
#     cells entirely (only the calculated Duration/columns remain, driven by
#     their formulas referencing the still-empty B27/C27). ---
$ws.Range("A27:C27").Clear()

# --- Extend the worksheet table (ListObject) & AutoFilter to include the
#     two new rows. ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F27"))

# --- Reflect the user's final scroll position / selection on the sheet ---
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
$ws.Range("D26").Select() | Out-Null
